# Applies the coinranking.com price-sheet refresh described in the commit
# "Updated cryptos list on Thu Oct  3 09:44:23 UTC 2024 with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'60.633.36"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -1.32%  '

$ws.Range('D3').Value = "'2.346.79"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -4.61%  '

$ws.Range('E4').Value = '  -0.02%  '

$ws.Range('D5').Value = "'542.31"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.62%  '

$ws.Range('D6').Value = "'137.06"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -6.56%  '

$ws.Range('E7').Value = '  -0.02%  '

$ws.Range('D8').Value = "'0.519"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -12.45%  '

$ws.Range('D9').Value = "'2.349.28"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -4.39%  '

$ws.Range('E10').Value = '  -3.62%  '

$ws.Range('E11').Value = '  -0.14%  '

$ws.Range('E12').Value = '  -4.56%  '

$ws.Range('E13').Value = '  -4.44%  '

$ws.Range('D14').Value = "'24.64"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -5.78%  '

$ws.Range('D15').Value = "'2.772.54"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -4.56%  '

$ws.Range('D16').Value = "'60.436.81"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.51%  '

$ws.Range('E17').Value = '  -4.06%  '

$ws.Range('D18').Value = "'2.347.28"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -4.67%  '

$ws.Range('E19').Value = '  -4.85%  '

$ws.Range('D20').Value = "'4.05"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -3.20%  '

$ws.Range('D21').Value = "'313.11"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.81%  '

$ws.Range('D22').Value = "'6.52"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -8.90%  '

$ws.Range('E23').Value = '  -0.17%  '

$ws.Range('D24').Value = "'1.87"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.45%  '

$ws.Range('D25').Value = "'62.83"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -2.20%  '

$ws.Range('D26').Value = "'8.17"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +6.38%  '

$ws.Range('D27').Value = "'0.998"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.20%  '

$ws.Range('D28').Value = "'2.463.60"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -5.35%  '

$ws.Range('D29').Value = "'7.93"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -4.25%  '

$ws.Range('D30').Value = "'0.0₃0883"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -10.72%  '

$ws.Range('D31').Value = "'499.62"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -10.92%  '

$ws.Range('E32').Value = '  -7.11%  '

$ws.Range('E33').Value = '  -3.04%  '

$ws.Range('E34').Value = '  -6.76%  '

$ws.Range('E35').Value = '  -4.03%  '

$ws.Range('E36').Value = '  -0.02%  '

$ws.Range('D37').Value = "'4.51"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -7.09%  '

$ws.Range('B38').Value = 'PolygonEcosystemToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D38').Value = "'0.370"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -2.51%  '

$ws.Range('B39').Value = 'EthereumClassic'
$ws.Range('C39').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D39').Value = "'18.28"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.72%  '

$ws.Range('D40').Value = "'5.21"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -11.54%  '

$ws.Range('D41').Value = "'1.76"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.98%  '

$ws.Range('B42').Value = 'Monero'
$ws.Range('C42').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D42').Value = "'138.28"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.73%  '

$ws.Range('B43').Value = 'USDe'
$ws.Range('C43').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D43').Value = "'1.00"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.02%  '

$ws.Range('D44').Value = "'40.02"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.32%  '

$ws.Range('D45').Value = "'138.58"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -5.49%  '

$ws.Range('D46').Value = "'3.51"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -2.71%  '

$ws.Range('E47').Value = '  -13.90%  '

$ws.Range('E48').Value = '  -5.07%  '

$ws.Range('D49').Value = "'19.41"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -10.18%  '

$ws.Range('E50').Value = '  -4.78%  '

$ws.Range('E51').Value = '  -5.18%  '
